$d = $word.ActiveDocument

# Move to the very end of the document content and add two new paragraphs:
#  - one empty paragraph (matching the existing trailing empty paragraph's style)
#  - one paragraph containing the text "test"
$endRange = $d.Content
$endRange.Collapse(0)            # wdCollapseEnd
$endRange.InsertParagraphAfter()
$endRange.Collapse(0)
$endRange.InsertAfter("test")
